# Update the "想去人数" (want-to-go count) and occasionally "最低票价" (min price)
# figures on the 展览 (Exhibition) and 全部类型 (All types) sheets, mirroring a
# fresh scrape (output generated at 456a3b4). 演出 (Performance) and
# 本地生活 (Local life) sheets are untouched by this refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 15160
$ws.Range("F3").Value = 19495
$ws.Range("G3").Value = 75
$ws.Range("F5").Value = 160
$ws.Range("F14").Value = 205
$ws.Range("F15").Value = 247
$ws.Range("F17").Value = 1511
$ws.Range("F20").Value = 111
$ws.Range("F21").Value = 247
$ws.Range("F22").Value = 8163
$ws.Range("F23").Value = 992
$ws.Range("F24").Value = 41
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 66
$ws.Range("F27").Value = 1269
$ws.Range("F28").Value = 10
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 6522
$ws.Range("F32").Value = 129
$ws.Range("F33").Value = 78
$ws.Range("F34").Value = 183
$ws.Range("F36").Value = 304
$ws.Range("F37").Value = 5551
$ws.Range("F38").Value = 1015
$ws.Range("F39").Value = 27
$ws.Range("F40").Value = 32
$ws.Range("F41").Value = 61

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 15160
$ws.Range("F3").Value = 19495
$ws.Range("G3").Value = 75
$ws.Range("F5").Value = 160
$ws.Range("F14").Value = 205
$ws.Range("F15").Value = 247
$ws.Range("F17").Value = 1511
$ws.Range("F21").Value = 111
$ws.Range("F22").Value = 247
$ws.Range("F23").Value = 8163
$ws.Range("F24").Value = 992
$ws.Range("F25").Value = 41
$ws.Range("F26").Value = 11
$ws.Range("F27").Value = 66
$ws.Range("F28").Value = 1269
$ws.Range("F29").Value = 10
$ws.Range("F30").Value = 13
$ws.Range("F31").Value = 20
$ws.Range("F34").Value = 6522
$ws.Range("F35").Value = 129
$ws.Range("F36").Value = 78
$ws.Range("F37").Value = 183
$ws.Range("F39").Value = 304
$ws.Range("F40").Value = 5551
$ws.Range("F41").Value = 1015
$ws.Range("F42").Value = 27
$ws.Range("F43").Value = 32
$ws.Range("F44").Value = 61
